$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns P1 and Q1, copying the formatting from O1 (bold, bordered, centered)
$ws.Range("O1").Copy($ws.Range("P1:Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For each data row (2..25):
#  - swap the values of I<->K and M<->O
#  - add new columns P and Q with value 2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # I
    $kVal = $ws.Cells.Item($r, 11).Value2  # K
    $mVal = $ws.Cells.Item($r, 13).Value2  # M
    $oVal = $ws.Cells.Item($r, 15).Value2  # O

    $ws.Cells.Item($r, 9).Value = $kVal    # I = old K
    $ws.Cells.Item($r, 11).Value = $iVal   # K = old I
    $ws.Cells.Item($r, 13).Value = $oVal   # M = old O
    $ws.Cells.Item($r, 15).Value = $mVal   # O = old M

    $ws.Cells.Item($r, 16).Value = 2       # P
    $ws.Cells.Item($r, 17).Value = 2       # Q
}
